# Lattice multiplication exercises worksheet: regenerate the 15 practice
# problems (5 rows x 3 columns) in the single table with a new set of
# multiplicand/multiplier pairs. Each cell holds one <w:r> run whose text
# is split across 5 lines by <w:br/> line breaks:
#   "AA x BB"
#   "  C    D"
#   "  ----"
#   "E|    |"
#   "F|    |"
# We rebuild each cell's Range.Text wholesale (using vertical-tab, chr(11),
# as the line-break character Word uses for <w:br/> inside Range.Text) so
# the existing run formatting (sz=32) is preserved.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$brk = [char]11

# row, col, [5 lines of new cell content]
$newCells = @(
    @(1, 1, @("61 x 35", "  3    5", "  ----", "6|    |", "1|    |")),
    @(1, 2, @("14 x 54", "  5    4", "  ----", "1|    |", "4|    |")),
    @(1, 3, @("13 x 95", "  9    5", "  ----", "1|    |", "3|    |")),
    @(2, 1, @("16 x 41", "  4    1", "  ----", "1|    |", "6|    |")),
    @(2, 2, @("20 x 22", "  2    2", "  ----", "2|    |", "0|    |")),
    @(2, 3, @("56 x 13", "  1    3", "  ----", "5|    |", "6|    |")),
    @(3, 1, @("65 x 92", "  9    2", "  ----", "6|    |", "5|    |")),
    @(3, 2, @("63 x 53", "  5    3", "  ----", "6|    |", "3|    |")),
    @(3, 3, @("78 x 17", "  1    7", "  ----", "7|    |", "8|    |")),
    @(4, 1, @("49 x 83", "  8    3", "  ----", "4|    |", "9|    |")),
    @(4, 2, @("78 x 14", "  1    4", "  ----", "7|    |", "8|    |")),
    @(4, 3, @("18 x 69", "  6    9", "  ----", "1|    |", "8|    |")),
    @(5, 1, @("76 x 80", "  8    0", "  ----", "7|    |", "6|    |")),
    @(5, 2, @("17 x 21", "  2    1", "  ----", "1|    |", "7|    |")),
    @(5, 3, @("57 x 64", "  6    4", "  ----", "5|    |", "7|    |"))
)

foreach ($entry in $newCells) {
    $row = $entry[0]
    $col = $entry[1]
    $lines = $entry[2]
    $text = [string]::Join($brk, $lines)
    $t.Cell($row, $col).Range.Text = $text
}

Write-Output "Updated $($newCells.Count) cells in table 1"
